$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1): PriceChange (X1), UpDown (Y1) ---
$ws.Cells.Item(1, 24).Value = "PriceChange"
$ws.Cells.Item(1, 25).Value = "UpDown"

# --- Overwrite existing row 2 with new data values ---
$ws.Cells.Item(2, 1).Value = 42633.878900462965
$ws.Cells.Item(2, 2).Value = 11
$ws.Cells.Item(2, 3).Value = "Buy"
$ws.Cells.Item(2, 4).Value = 34
$ws.Cells.Item(2, 5).Value = 10620
$ws.Cells.Item(2, 6).Value = 1266
$ws.Cells.Item(2, 7).Value = 62
$ws.Cells.Item(2, 8).Value = 36
$ws.Cells.Item(2, 9).Value = 94
$ws.Cells.Item(2, 10).Value = 5
$ws.Cells.Item(2, 11).Value = 6892
$ws.Cells.Item(2, 12).Value = 139
$ws.Cells.Item(2, 13).Value = 81
$ws.Cells.Item(2, 14).Value = 16
$ws.Cells.Item(2, 15).Value = 1
$ws.Cells.Item(2, 16).Value = "Bag"
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = 1.76
$ws.Cells.Item(2, 19).Value = 0.111
$ws.Cells.Item(2, 20).Value = 5.45
$ws.Cells.Item(2, 21).Value = 4.84
$ws.Cells.Item(2, 22).Value = 2.2799999999999998
$ws.Cells.Item(2, 23).Value = 0
$ws.Cells.Item(2, 24).Value = -1.6100000000000136
$ws.Cells.Item(2, 25).Value = "Down"

# --- New row 3 (copy formats from row 2 first so date/percent styles match) ---
$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("S2").Copy($ws.Range("S3"))

$ws.Cells.Item(3, 1).Value = 42633.880312499998
$ws.Cells.Item(3, 2).Value = 11
$ws.Cells.Item(3, 3).Value = "Buy"
$ws.Cells.Item(3, 4).Value = 44
$ws.Cells.Item(3, 5).Value = 10847
$ws.Cells.Item(3, 6).Value = 1310
$ws.Cells.Item(3, 7).Value = 67
$ws.Cells.Item(3, 8).Value = 32
$ws.Cells.Item(3, 9).Value = 100
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 7060
$ws.Cells.Item(3, 12).Value = 146
$ws.Cells.Item(3, 13).Value = 70
$ws.Cells.Item(3, 14).Value = 20
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = "Bag"
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 1.76
$ws.Cells.Item(3, 19).Value = 0.111
$ws.Cells.Item(3, 20).Value = 5.45
$ws.Cells.Item(3, 21).Value = 4.84
$ws.Cells.Item(3, 22).Value = 2.2799999999999998
$ws.Cells.Item(3, 23).Value = 0

# --- Column width tweak (column C); stored OOXML width = ColumnWidth + 0.8333333333333334 ---
# (closest reachable value to the target stored width of 7.42578125)
$ws.Columns("C").ColumnWidth = 6.666666666666667
